$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Runner")

# Update cell C2 from "no" to "yes"
$ws.Range("C2").Value = "yes"

# Update the active selection to C2 (as recorded in the saved view state)
$ws.Range("C2").Select()
